# Bugfix: anzahlRubine werden bei next level vom FMENU nicht resetet
# -> add a separate "Schlange X" tracking entry alongside the existing
#    "Schlange" entry (renamed to "Schlange Y") in the RGB value table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Schlange" row (row 7) to "Schlange Y"
$ws.Range("A7").Value = "Schlange Y"

# Insert a new row for "Schlange X" right above the old row 8 ("Spawn"),
# pushing Spawn/Stein/Wand/Weg down by one row
$ws.Rows("8:8").Insert()

# Populate the new row 8 with the "Schlange X" RGB data
$ws.Range("A8").Value = "Schlange X"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 250
$ws.Range("D8").Value = 0

# Match the saved selection state
$ws.Range("A8").Select()

# Re-apply the table sort (by Objekt, ascending) over the now-larger range
# so the sheet's remembered sort state covers A2:D12 / A2:A12
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A12"))
$ws.Sort.SetRange($ws.Range("A2:D12"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()
